$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Camote" needs to be added to the
# dataset. It belongs chronologically/logically at the top of this block
# of rows, so insert a fresh row at 240 and push the existing rows
# (240-245) down to (241-246).
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new observation's data.
$ws.Cells.Item(240, 1).Value = 10
$ws.Cells.Item(240, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(240, 3).Value = "La Araucanía"
$ws.Cells.Item(240, 4).Value = 45239
$ws.Cells.Item(240, 5).Value = 9
$ws.Cells.Item(240, 6).Value = 100114002
$ws.Cells.Item(240, 7).Value = "Camote"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 100
$ws.Cells.Item(240, 11).Value = 24000
$ws.Cells.Item(240, 12).Value = 24000
$ws.Cells.Item(240, 13).Value = 24000
$ws.Cells.Item(240, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(240, 15).Value = "Perú"
$ws.Cells.Item(240, 16).Value = 1333
$ws.Cells.Item(240, 17).Value = 18
$ws.Cells.Item(240, 18).Value = "Hortaliza"
